# Generate Report for Handoff
#
# The localization-status workbook tracks three "source" files
# (6104303d-...md, ffff5f0b0cd7-...md, ffffff7a5bff45-...md) per
# language sheet. This re-generates the report: the row that used to
# describe 6104303d-...md (previously "Handed back: in sync with en-US")
# is now re-issued ("Ready for handoff") with fresh handoff timestamps,
# and the remaining two rows shift up to take its old slot.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (File Name / zh-cn / de-de) ----
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value2 = "ffff5f0b0cd7-2e4d-4fac-8347-038a6aff8eeb.md"

$ws.Range("A3").Value2 = "ffffff7a5bff45-0784-4e83-b89b-0b0d1af02392.md"

$ws.Range("A4").Value2 = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "Ready for handoff"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value2 = "ffff5f0b0cd7-2e4d-4fac-8347-038a6aff8eeb.md"
$ws.Range("C2").Value2 = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.zh-cn.xlf"
$ws.Range("D2").Value2 = "2016-03-09 10:06:00"
$ws.Range("E2").Value2 = "acd3cc98-32a6-43cb-9bfd-62d79904db49.md"
$ws.Range("F2").Value2 = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.zh-cn.xlf"
$ws.Range("G2").Value2 = "2016-03-09 10:06:29"

$ws.Range("A3").Value2 = "ffffff7a5bff45-0784-4e83-b89b-0b0d1af02392.md"

$ws.Range("A4").Value2 = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "6104303d-37b8-4b92-8309-bfe68b998bc1.507280672b761ae06d2d43713a27199e770b6384.zh-cn.xlf"
$ws.Range("D4").Value2 = "2016-03-09 10:08:47"
$ws.Range("E4").Value2 = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws.Range("F4").Value2 = "6104303d-37b8-4b92-8309-bfe68b998bc1.507280672b761ae06d2d43713a27199e770b6384.zh-cn.xlf"
$ws.Range("G4").Value2 = "2016-03-09 10:08:23"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value2 = "ffff5f0b0cd7-2e4d-4fac-8347-038a6aff8eeb.md"
$ws.Range("C2").Value2 = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.de-de.xlf"
$ws.Range("D2").Value2 = "2016-03-09 10:06:06"
$ws.Range("E2").Value2 = "acd3cc98-32a6-43cb-9bfd-62d79904db49.md"
$ws.Range("F2").Value2 = "acd3cc98-32a6-43cb-9bfd-62d79904db49.aab57bf76b38a394b6610a9034d9b6ef5852519a.de-de.xlf"
$ws.Range("G2").Value2 = "2016-03-09 10:06:35"

$ws.Range("A3").Value2 = "ffffff7a5bff45-0784-4e83-b89b-0b0d1af02392.md"

$ws.Range("A4").Value2 = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "6104303d-37b8-4b92-8309-bfe68b998bc1.507280672b761ae06d2d43713a27199e770b6384.de-de.xlf"
$ws.Range("D4").Value2 = "2016-03-09 10:08:51"
$ws.Range("E4").Value2 = "6104303d-37b8-4b92-8309-bfe68b998bc1.md"
$ws.Range("F4").Value2 = "6104303d-37b8-4b92-8309-bfe68b998bc1.507280672b761ae06d2d43713a27199e770b6384.de-de.xlf"
$ws.Range("G4").Value2 = "2016-03-09 10:08:29"
